# Auto-generated edit script applying numeric corrections to Kujata_Profits sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 3124.1177
$ws.Range("I6").Value = 5656.6665
$ws.Range("K6").Value = 16969.9995
$ws.Range("M6").Value = -16857.9995

$ws.Range("H31").Value = 312
$ws.Range("I31").Value = 312
$ws.Range("K31").Value = 936
$ws.Range("M31").Value = -706

$ws.Range("H42").Value = 776.3333
$ws.Range("I42").Value = 936
$ws.Range("J42").Value = 616.6667
$ws.Range("K42").Value = 2808
$ws.Range("L42").Value = 1850.0001
$ws.Range("M42").Value = -2578
$ws.Range("N42").Value = -2310.0001

$ws.Range("H97").Value = 2187.4092
$ws.Range("J97").Value = 2187.4092
$ws.Range("L97").Value = 6562.2276
$ws.Range("N97").Value = -7554.2276

$ws.Range("H106").Value = 2166.8948
$ws.Range("I106").Value = 2145
$ws.Range("J106").Value = 2353
$ws.Range("K106").Value = 2145
$ws.Range("L106").Value = 2353
$ws.Range("M106").Value = -1514
$ws.Range("N106").Value = -3615

$ws.Range("H137").Value = 1401.2
$ws.Range("I137").Value = 1129.7142
$ws.Range("J137").Value = 1638.75
$ws.Range("K137").Value = 3389.1426
$ws.Range("L137").Value = 4916.25
$ws.Range("M137").Value = -839.1425999999997
$ws.Range("N137").Value = -10016.25

$ws.Range("H138").Value = 2040.2279
$ws.Range("I138").Value = 1449.7142
$ws.Range("J138").Value = 2167.4153
$ws.Range("K138").Value = 4349.142599999999
$ws.Range("L138").Value = 6502.2459
$ws.Range("M138").Value = 790.8574000000008
$ws.Range("N138").Value = -16782.2459

$ws.Range("H141").Value = 7368.8887
$ws.Range("I141").Value = 8395.357
$ws.Range("J141").Value = 3776.25
$ws.Range("K141").Value = 25186.071
$ws.Range("L141").Value = 11328.75
$ws.Range("M141").Value = -20006.071
$ws.Range("N141").Value = -21688.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4738.971
$ws.Range("I32").Value = 4738.971
$ws.Range("K32").Value = 4738.971
$ws.Range("M32").Value = -4451.971

$ws.Range("H49").Value = 10036
$ws.Range("I49").Value = 10036
$ws.Range("J49").Value = 0
$ws.Range("K49").Value = 10036
$ws.Range("L49").Value = 0
$ws.Range("M49").Value = -9776
$ws.Range("N49").ClearContents()

$ws.Range("H132").Value = 6438.5
$ws.Range("I132").Value = 8427.5
$ws.Range("J132").Value = 4449.5
$ws.Range("K132").Value = 25282.5
$ws.Range("L132").Value = 13348.5
$ws.Range("M132").Value = -22752.5
$ws.Range("N132").Value = -18408.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1345.8889
$ws.Range("I20").Value = 1170.5385
$ws.Range("J20").Value = 1801.8
$ws.Range("K20").Value = 1170.5385
$ws.Range("L20").Value = 1801.8
$ws.Range("M20").Value = -923.5385000000001
$ws.Range("N20").Value = -2295.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 6660046.5
$ws.Range("I6").Value = 6660046.5
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 6660046.5
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -6659933.5
$ws.Range("N6").ClearContents()

$ws.Range("H16").Value = 111112264
$ws.Range("I16").Value = 166667710
$ws.Range("J16").Value = 1366.6666
$ws.Range("K16").Value = 166667710
$ws.Range("L16").Value = 1366.6666
$ws.Range("M16").Value = -166667423
$ws.Range("N16").Value = -1940.6666

$ws.Range("H59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").ClearContents()

$ws.Range("H113").Value = 111112264
$ws.Range("I113").Value = 166667710
$ws.Range("J113").Value = 1366.6666
$ws.Range("K113").Value = 166667710
$ws.Range("L113").Value = 1366.6666
$ws.Range("M113").Value = -166665540
$ws.Range("N113").Value = -5706.6666

$ws.Range("H132").Value = 8436.177
$ws.Range("I132").Value = 11700.3
$ws.Range("J132").Value = 3773.1428
$ws.Range("K132").Value = 35100.89999999999
$ws.Range("L132").Value = 11319.4284
$ws.Range("M132").Value = -32570.89999999999
$ws.Range("N132").Value = -16379.4284

$ws.Range("H141").Value = 34000
$ws.Range("I141").Value = 34000
$ws.Range("J141").Value = 34000
$ws.Range("K141").Value = 34000
$ws.Range("L141").Value = 34000
$ws.Range("M141").Value = -28820
$ws.Range("N141").Value = -44360

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1440.919
$ws.Range("I68").Value = 770.55554
$ws.Range("J68").Value = 1656.3928
$ws.Range("K68").Value = 2311.66662
$ws.Range("L68").Value = 4969.178400000001
$ws.Range("M68").Value = -1500.66662
$ws.Range("N68").Value = -6591.178400000001

$ws.Range("H71").Value = 1440.919
$ws.Range("I71").Value = 770.55554
$ws.Range("J71").Value = 1656.3928
$ws.Range("K71").Value = 6934.99986
$ws.Range("L71").Value = 14907.5352
$ws.Range("M71").Value = -2878.99986
$ws.Range("N71").Value = -23019.5352

$ws.Range("H104").Value = 5747.3335
$ws.Range("I104").Value = 3113
$ws.Range("J104").Value = 6500
$ws.Range("K104").Value = 9339
$ws.Range("L104").Value = 19500
$ws.Range("M104").Value = -6718
$ws.Range("N104").Value = -24742

$ws.Range("H122").Value = 733.4666999999999
$ws.Range("I122").Value = 480.8889
$ws.Range("J122").Value = 1112.3334
$ws.Range("K122").Value = 4328.0001
$ws.Range("L122").Value = 10011.0006
$ws.Range("M122").Value = -1878.0001
$ws.Range("N122").Value = -14911.0006

$ws.Range("H123").Value = 0
$ws.Range("I123").Value = 0
$ws.Range("K123").Value = 0
$ws.Range("M123").ClearContents()

$ws.Range("H131").Value = 16950452
$ws.Range("I131").Value = 111111900
$ws.Range("J131").Value = 1390.24
$ws.Range("K131").Value = 333335700
$ws.Range("L131").Value = 4170.72
$ws.Range("M131").Value = -333330660
$ws.Range("N131").Value = -14250.72

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2366.5
$ws.Range("I126").Value = 1733
$ws.Range("J126").Value = 3000
$ws.Range("K126").Value = 5199
$ws.Range("L126").Value = 9000
$ws.Range("M126").Value = -2729
$ws.Range("N126").Value = -13940

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2185.2856
$ws.Range("I7").Value = 1965.3334
$ws.Range("K7").Value = 1965.3334
$ws.Range("M7").Value = -1853.3334

$ws.Range("H46").Value = 5108.5
$ws.Range("I46").Value = 350.1111
$ws.Range("K46").Value = 350.1111
$ws.Range("M46").Value = -162.1111

$ws.Range("H47").Value = 5532.5
$ws.Range("I47").Value = 3000
$ws.Range("J47").Value = 8065
$ws.Range("K47").Value = 3000
$ws.Range("L47").Value = 8065
$ws.Range("M47").Value = -2510
$ws.Range("N47").Value = -9045

$ws.Range("H52").Value = 5532.5
$ws.Range("I52").Value = 3000
$ws.Range("J52").Value = 8065
$ws.Range("K52").Value = 3000
$ws.Range("L52").Value = 8065
$ws.Range("M52").Value = -2767
$ws.Range("N52").Value = -8531

$ws.Range("H126").Value = 2185.2856
$ws.Range("I126").Value = 1965.3334
$ws.Range("K126").Value = 5896.0002
$ws.Range("M126").Value = -3426.0002

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H108").Value = 11541.333
$ws.Range("J108").Value = 11541.333
$ws.Range("L108").Value = 11541.333
$ws.Range("N108").Value = -19221.333
